# Append the 02/27/2026 Kaspa buy as a new row at the bottom of the sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A stores the date as literal text (e.g. "02/20/2026" in the row
# above), not a real Excel date serial, so force the cell to Text before
# assigning the value to keep it from being auto-converted to a date.
$ws.Range("A29").NumberFormat = "@"
$ws.Range("A29").Value = "02/27/2026"

$ws.Range("B29").Value = 1557.662000000004
$ws.Range("C29").Value = 0.03177839608336075
$ws.Range("D29").Value = 50
